$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price cells to remain text (matches source t="inlineStr")
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated values
$ws.Range('D2').Value = '63.179.17'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.558.88'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '583.65'
$ws.Range('E5').Value = '  +2.26%  '
$ws.Range('D6').Value = '147.59'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('E9').Value = '  +3.24%  '
$ws.Range('D10').Value = '5.61'
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '0.356'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '27.59'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '3.016.67'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '63.101.46'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('E16').Value = '  +4.34%  '
$ws.Range('D17').Value = '2.562.94'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '11.36'
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('D19').Value = '341.12'
$ws.Range('E19').Value = '  +1.79%  '
$ws.Range('D20').Value = '4.42'
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('D21').Value = '6.80'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('D24').Value = '2.685.16'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('D27').Value = '8.07'
$ws.Range('E27').Value = '  +12.35%  '
$ws.Range('D28').Value = '8.54'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +8.32%  '
$ws.Range('D32').Value = '0.0₃0823'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').Value = '443.27'
$ws.Range('E33').Value = '  +6.34%  '
$ws.Range('D34').Value = '176.81'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').Value = '1.62'
$ws.Range('E35').Value = '  +1.97%  '
$ws.Range('E36').Value = '  +2.48%  '
$ws.Range('D37').Value = '19.26'
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').Value = '150.67'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '3.83'
$ws.Range('E43').Value = '  +1.87%  '
$ws.Range('D44').Value = '21.07'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('D45').Value = '0.0552'
$ws.Range('E45').Value = '  +6.26%  '
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').Value = '0.0975'
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D49').Value = '18.41'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  -3.25%  '
$ws.Range('E51').Value = '  -0.31%  '

# Strip the temporary text-format styling so cell style index matches the original (unstyled)
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
